$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 and IF in columns I and J, matching the style of the
# existing header row (row 1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF) for rows 2-21
$data = @(
    @(1,1),
    @(7,7),
    @(6,7),
    @(1,2),
    @(6,7),
    @(8,9),
    @(6,6),
    @(6,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(6,7),
    @(4,5),
    @(1,3),
    @(6,6),
    @(6,6),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
